# Update average_county_temperature (I), worst_ashp_cop (N) and
# best_ashp_cop (O) for the affected facility rows using refreshed
# NOAA county-temperature data (merged dataset update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 19.79629629629628
$ws.Range("N3").Value = 1.911855479578636
$ws.Range("O3").Value = 2.09608909874769
$ws.Range("I13").Value = 21.28240740740739
$ws.Range("I14").Value = 21.28240740740739
$ws.Range("N14").Value = 1.940636870984383
$ws.Range("O14").Value = 2.131200751448103
$ws.Range("I17").Value = 19.30324074074072
$ws.Range("N17").Value = 1.902494195535734
$ws.Range("O17").Value = 2.084694111942012
$ws.Range("I18").Value = 15.74228395061728
$ws.Range("N18").Value = 1.837513876759573
$ws.Range("O18").Value = 2.005936573945218
$ws.Range("I21").Value = 20.22222222222222
$ws.Range("N21").Value = 1.920016703786191
$ws.Range("O21").Value = 2.106033415841584
$ws.Range("I22").Value = 13.75752314814816
$ws.Range("N22").Value = 1.803186500133452
$ws.Range("O22").Value = 1.964569140204562
$ws.Range("I25").Value = 12.93898809523811
$ws.Range("N25").Value = 1.789400236291612
$ws.Range("O25").Value = 1.948001533154466
$ws.Range("I26").Value = 12.93898809523811
$ws.Range("I31").Value = 5.486111111111112
$ws.Range("N31").Value = 1.672941176470588
$ws.Range("O31").Value = 1.809089700996678
$ws.Range("I32").Value = 21.28240740740739
$ws.Range("N32").Value = 1.940636870984383
$ws.Range("O32").Value = 2.131200751448103
$ws.Range("I36").Value = 1.791666666666668
$ws.Range("N36").Value = 1.620655622136059
$ws.Range("O36").Value = 1.747323835194455
$ws.Range("I37").Value = 5.486111111111112
$ws.Range("N37").Value = 1.672941176470588
$ws.Range("O37").Value = 1.809089700996678
$ws.Range("I41").Value = 19.30324074074072
$ws.Range("N41").Value = 1.902494195535734
$ws.Range("O41").Value = 2.084694111942012
$ws.Range("I43").Value = 12.67039049919483
$ws.Range("N43").Value = 1.784922174701128
$ws.Range("O43").Value = 1.942625691911729
$ws.Range("I49").Value = 21.28240740740739
$ws.Range("N49").Value = 1.940636870984383
$ws.Range("O49").Value = 2.131200751448103
$ws.Range("I52").Value = 13.62268518518517
$ws.Range("N52").Value = 1.80090088129692
$ws.Range("O52").Value = 1.961820583643568
$ws.Range("I55").Value = 12.67039049919483
$ws.Range("N55").Value = 1.784922174701128
$ws.Range("O55").Value = 1.942625691911729
$ws.Range("I57").Value = 12.67039049919483
$ws.Range("N57").Value = 1.784922174701128
$ws.Range("O57").Value = 1.942625691911729
$ws.Range("I60").Value = 19.65277777777778
$ws.Range("N60").Value = 1.909121107266436
$ws.Range("O60").Value = 2.092759415833974
$ws.Range("I68").Value = 19.30324074074072
$ws.Range("N68").Value = 1.902494195535734
$ws.Range("O68").Value = 2.084694111942012
$ws.Range("I77").Value = 13.76976495726495
$ws.Range("N77").Value = 1.803394296576035
$ws.Range("O77").Value = 1.964819060413116
$ws.Range("I80").Value = 21.28240740740739
$ws.Range("N80").Value = 1.940636870984383
$ws.Range("O80").Value = 2.131200751448103
$ws.Range("I81").Value = 13.75752314814816
$ws.Range("N81").Value = 1.803186500133452
$ws.Range("O81").Value = 1.964569140204562
$ws.Range("I82").Value = 13.75752314814816
